$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, reusing the same format as the other header cells
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add 0 values in H2:H7 (plain numeric cells, no special style)
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}

$excel.CutCopyMode = 0
